$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = New-Object 'double[,]' 24,11

$data[0,0] = 2.826683192430096
$data[0,1] = 0.3648183455056255
$data[0,2] = 0.009022407013876688
$data[0,3] = 0.05043111356057539
$data[0,4] = 4.157605467521449
$data[0,5] = 0
$data[0,6] = 0.07973214163530429
$data[0,7] = 2.446716918797193
$data[0,8] = 0.1247199570505479
$data[0,9] = 0
$data[0,10] = 0.4654450911102828
$data[1,0] = 2.736422786721732
$data[1,1] = 0.3374530303529468
$data[1,2] = 0.008450282533296161
$data[1,3] = 0.05050730774578205
$data[1,4] = 4.128671466143004
$data[1,5] = 0
$data[1,6] = 0.07973214163530429
$data[1,7] = 2.436935429613129
$data[1,8] = 0.1252668101822341
$data[1,9] = 0
$data[1,10] = 0.4601784076901652
$data[2,0] = 2.682727031003594
$data[2,1] = 0.3208750531883595
$data[2,2] = 0.008096173152669195
$data[2,3] = 0.05055825357180782
$data[2,4] = 4.112946417055952
$data[2,5] = 0
$data[2,6] = 0.07973214163530429
$data[2,7] = 2.432049083497702
$data[2,8] = 0.1256277414360021
$data[2,9] = 0
$data[2,10] = 0.457184764822486
$data[3,0] = 2.661278921277471
$data[3,1] = 0.3141752933427426
$data[3,2] = 0.007951108721673705
$data[3,3] = 0.05058006263722614
$data[3,4] = 4.107049310230821
$data[3,5] = 0
$data[3,6] = 0.07973214163530429
$data[3,7] = 2.4303381143487
$data[3,8] = 0.1257811514933511
$data[3,9] = 0
$data[3,10] = 0.4560252048456164
$data[4,0] = 2.657743644695699
$data[4,1] = 0.3130661641413042
$data[4,2] = 0.007926973338296506
$data[4,3] = 0.05058374737421739
$data[4,4] = 4.106100899938554
$data[4,5] = 0
$data[4,6] = 0.07973214163530429
$data[4,7] = 2.430070899375011
$data[4,8] = 0.1258070073103639
$data[4,9] = 0
$data[4,10] = 0.4558363071264608
$data[5,0] = 2.682436020064245
$data[5,1] = 0.3207844722414279
$data[5,2] = 0.008094219916408463
$data[5,3] = 0.05055854344985866
$data[5,4] = 4.112864820535606
$data[5,5] = 0
$data[5,6] = 0.07973214163530429
$data[5,7] = 2.43202487571115
$data[5,8] = 0.1256297847564181
$data[5,9] = 0
$data[5,10] = 0.4571688821454671
$data[6,0] = 2.795203256876675
$data[6,1] = 0.3553358682499095
$data[6,2] = 0.008825692755543457
$data[6,3] = 0.05045652236715681
$data[6,4] = 4.14720420287432
$data[6,5] = 0
$data[6,6] = 0.07973214163530429
$data[6,7] = 2.443111137013076
$data[6,8] = 0.1249032902062606
$data[6,9] = 0
$data[6,10] = 0.4635792929911844
$data[7,0] = 3.03005399005832
$data[7,1] = 0.4249002359609904
$data[7,2] = 0.01024010721229729
$data[7,3] = 0.05028941902571737
$data[7,4] = 4.230838823666346
$data[7,5] = 0
$data[7,6] = 0.07973214163530429
$data[7,7] = 2.473794460819761
$data[7,8] = 0.1236782477282627
$data[7,9] = 0
$data[7,10] = 0.478057153462089
$data[8,0] = 3.211029415870769
$data[8,1] = 0.4771574687177349
$data[8,2] = 0.01127045713518626
$data[8,3] = 0.05018665766128072
$data[8,4] = 4.302375552408307
$data[8,5] = 0
$data[8,6] = 0.07973214163530429
$data[8,7] = 2.501878032762775
$data[8,8] = 0.1228998502988183
$data[8,9] = 0
$data[8,10] = 0.4898613003191628
$data[9,0] = 3.295208059933202
$data[9,1] = 0.5011903157214874
$data[9,2] = 0.01173800890599352
$data[9,3] = 0.0501442375544463
$data[9,4] = 4.337145307364466
$data[9,5] = 0
$data[9,6] = 0.07973214163530429
$data[9,7] = 2.515876609817639
$data[9,8] = 0.122572142332988
$data[9,9] = 0
$data[9,10] = 0.4954859969533487
$data[10,0] = 3.32735167314712
$data[10,1] = 0.5103291676626895
$data[10,2] = 0.01191495070814952
$data[10,3] = 0.05012879513045387
$data[10,4] = 4.350634647705874
$data[10,5] = 0
$data[10,6] = 0.07973214163530429
$data[10,7] = 2.52135491836404
$data[10,8] = 0.1224518426528665
$data[10,9] = 0
$data[10,10] = 0.4976526515258968
$data[11,0] = 3.320417085226381
$data[11,1] = 0.5083592460607065
$data[11,2] = 0.01187684710297887
$data[11,3] = 0.05013209332233737
$data[11,4] = 4.347715080630792
$data[11,5] = 0
$data[11,6] = 0.07973214163530429
$data[11,7] = 2.52016715554231
$data[11,8] = 0.122477582478874
$data[11,9] = 0
$data[11,10] = 0.4971843904977788
$data[12,0] = 3.297847179447501
$data[12,1] = 0.5019414066046579
$data[12,2] = 0.01175256789003853
$data[12,3] = 0.05014295465366292
$data[12,4] = 4.338248599914436
$data[12,5] = 0
$data[12,6] = 0.07973214163530429
$data[12,7] = 2.516323750798534
$data[12,8] = 0.1225621691147705
$data[12,9] = 0
$data[12,10] = 0.4956635130473614
$data[13,0] = 3.284057251091326
$data[13,1] = 0.4980152799802795
$data[13,2] = 0.01167643066424162
$data[13,3] = 0.05014968839639733
$data[13,4] = 4.332492218495076
$data[13,5] = 0
$data[13,6] = 0.07973214163530429
$data[13,7] = 2.513992694172885
$data[13,8] = 0.1226144752979668
$data[13,9] = 0
$data[13,10] = 0.494736713078737
$data[14,0] = 3.205565466644885
$data[14,1] = 0.4755921670246153
$data[14,2] = 0.01123988304264145
$data[14,3] = 0.05018951688886414
$data[14,4] = 4.3001482999004
$data[14,5] = 0
$data[14,6] = 0.07973214163530429
$data[14,7] = 2.500987931072103
$data[14,8] = 0.1229217979256561
$data[14,9] = 0
$data[14,10] = 0.4894988467749073
$data[15,0] = 3.157888101739502
$data[15,1] = 0.4619035507006402
$data[15,2] = 0.01097181851940121
$data[15,3] = 0.05021505777589957
$data[15,4] = 4.280878597272391
$data[15,5] = 0
$data[15,6] = 0.07973214163530429
$data[15,7] = 2.493324255758509
$data[15,8] = 0.1231170902942615
$data[15,9] = 0
$data[15,10] = 0.4863509098895946
$data[16,0] = 3.130639618979217
$data[16,1] = 0.4540547265779651
$data[16,2] = 0.01081752232324718
$data[16,3] = 0.05023015550986099
$data[16,4] = 4.270004670293247
$data[16,5] = 0
$data[16,6] = 0.07973214163530429
$data[16,7] = 2.489031359859126
$data[16,8] = 0.1232319010555862
$data[16,9] = 0
$data[16,10] = 0.4845642921301732
$data[17,0] = 3.121443663655612
$data[17,1] = 0.4514014393771504
$data[17,2] = 0.01076525922535865
$data[17,3] = 0.05023533732760382
$data[17,4] = 4.266358850947427
$data[17,5] = 0
$data[17,6] = 0.07973214163530429
$data[17,7] = 2.487597574498267
$data[17,8] = 0.1232712005489187
$data[17,9] = 0
$data[17,10] = 0.4839634934606636
$data[18,0] = 3.16294539701704
$data[18,1] = 0.463358186445987
$data[18,2] = 0.01100036568513829
$data[18,3] = 0.05021229675946842
$data[18,4] = 4.28290818993608
$data[18,5] = 0
$data[18,6] = 0.07973214163530429
$data[18,7] = 2.494128149624217
$data[18,8] = 0.1230960440186184
$data[18,9] = 0
$data[18,10] = 0.4866835297257808
$data[19,0] = 3.304469256256823
$data[19,1] = 0.5038254415258621
$data[19,2] = 0.01178907423483722
$data[19,3] = 0.0501397475671902
$data[19,4] = 4.34102035429305
$data[19,5] = 0
$data[19,6] = 0.07973214163530429
$data[19,7] = 2.517447827015786
$data[19,8] = 0.1225372209288373
$data[19,9] = 0
$data[19,10] = 0.496109235096398
$data[20,0] = 3.398519843450401
$data[20,1] = 0.5304957142214448
$data[20,2] = 0.01230392696530558
$data[20,3] = 0.05009595260093613
$data[20,4] = 4.380882655891156
$data[20,5] = 0
$data[20,6] = 0.07973214163530429
$data[20,7] = 2.533723008122294
$data[20,8] = 0.122194125747054
$data[20,9] = 0
$data[20,10] = 0.5024834424987716
$data[21,0] = 3.348180611549708
$data[21,1] = 0.5162407074031989
$data[21,2] = 0.01202917764263844
$data[21,3] = 0.05011899586912338
$data[21,4] = 4.359434298955534
$data[21,5] = 0
$data[21,6] = 0.07973214163530429
$data[21,7] = 2.524941503606982
$data[21,8] = 0.1223752168633379
$data[21,9] = 0
$data[21,10] = 0.4990618139456586
$data[22,0] = 3.160658489683897
$data[22,1] = 0.4627004804384569
$data[22,2] = 0.01098746008345586
$data[22,3] = 0.05021354372555259
$data[22,4] = 4.281989974366326
$data[22,5] = 0
$data[22,6] = 0.07973214163530429
$data[22,7] = 2.493764357182002
$data[22,8] = 0.1231055511456507
$data[22,9] = 0
$data[22,10] = 0.4865330801220864
$data[23,0] = 2.965045411251083
$data[23,1] = 0.405882907746502
$data[23,2] = 0.009859264736508777
$data[23,3] = 0.05033110525095957
$data[23,4] = 4.20645286642096
$data[23,5] = 0
$data[23,6] = 0.07973214163530429
$data[23,7] = 2.464527273448056
$data[23,8] = 0.1239882895614226
$data[23,9] = 0
$data[23,10] = 0.473935881096466

$ws.Range("B2:L25").Value2 = $data
Write-Host "Applied updates to B2:L25"